$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6 (current row 6 = "step_3") to make room for
# the new "step_2_motion" entry, shifting everything below down by one.
$ws.Rows("6:6").Insert()

# Fill in the new row 6 with the motion-step translations.
$ws.Range("A6").Value = "step_2_motion"
$ws.Range("B6").Value = "Shake your phone!"
$ws.Range("C6").Value = "Schüttel dein Handy!"

# Drop the trailing colon from the mode labels (now on rows 8 and 10 after
# the insert).
$ws.Range("B8").Value = "Implied Square"
$ws.Range("C8").Value = "Angedeutetes Viereck"

$ws.Range("C10").Value = "Neon-Form"
$ws.Range("B10").Value = "Neon Shape"

# Match the saved selection from the author's session.
$ws.Range("B10").Select()
